$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.334.98"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "2.931.52"
$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'594.60"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").Value = "'143.47"
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -1.27%  "

$ws.Range("D9").Value = "'6.94"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  -1.88%  "

$ws.Range("D11").Value = "'0.437"
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  -1.14%  "

$ws.Range("D13").Value = "'33.22"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "3.415.90"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").Value = "61.336.58"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "2.927.74"
$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("D18").Value = "'6.64"
$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").Value = "'433.50"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").Value = "'13.53"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22").Value = "'7.05"
$ws.Range("E22").Value = "  -0.30%  "

$ws.Range("D23").Value = "'81.47"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'10.84"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("E25").Value = "  -2.14%  "

$ws.Range("D26").Value = "'11.71"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -3.72%  "

$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").Value = "'6.88"
$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("D31").Value = "'26.63"
$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("E33").Value = "  +0.08%  "

$ws.Range("E34").Value = "  +1.86%  "

$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("E36").Value = "  -0.31%  "

$ws.Range("D37").Value = "'2.96"
$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").Value = "'8.49"
$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("D41").Value = "'42.10"
$ws.Range("E41").Value = "  +5.19%  "

$ws.Range("D42").Value = "'0.279"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").Value = "2.703.54"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("D45").Value = "'133.57"
$ws.Range("E45").Value = "  +1.52%  "

$ws.Range("D46").Value = "'363.08"
$ws.Range("E46").Value = "  -3.22%  "

$ws.Range("D48").Value = "'23.51"
$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("E49").Value = "  -1.34%  "

$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("E51").Value = "  -0.05%  "
